# Update the Correspond Handoff Datetime / Correspond Handback DateTime
# columns on the zh-cn and de-de report sheets to reflect the regenerated
# handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-18 04:05:35"
$wsZhCn.Range("E5").Value = "2016-03-18 04:05:35"
$wsZhCn.Range("H4").Value = "2016-03-18 04:05:53"
$wsZhCn.Range("H5").Value = "2016-03-18 04:05:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-18 04:05:38"
$wsDeDe.Range("E5").Value = "2016-03-18 04:05:38"
$wsDeDe.Range("H4").Value = "2016-03-18 04:05:58"
$wsDeDe.Range("H5").Value = "2016-03-18 04:05:58"
